# ------------------------------------------------------------------
# Reproduces the commit "Adição de macro EXCEL, xlwings e histórico."
#   - Gastos: wipes the sample rows, adds a "Pendente" column (G)
#   - Salas: brand-new sheet + table (ID_Sala / NomeSala / ValorTotalGastos)
#   - Adicionar_Gastos: ID_Sala seed values become 1 / 2, rest of the
#     sample rows are wiped
#   - Home: brand-new, blank, becomes the active tab
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Gastos: clear the demo data but keep the formatting, then add
#    a 7th ("Pendente") column to the table.
# ---------------------------------------------------------------
$wsGastos = $wb.Worksheets.Item("Gastos")
$wsGastos.Range("A2:F18").ClearContents()

$loGastos = $wsGastos.ListObjects.Item(1)
$loGastos.ListColumns.Add() | Out-Null

$wsGastos.Range("G1").Value = "Pendente"
$wsGastos.Range("G1").Font.Bold = $true
$wsGastos.Range("G1").Font.Color = 0xFFFFFF
$wsGastos.Range("G1").Interior.Color = 0x000000
$wsGastos.Range("G1").HorizontalAlignment = -4108
$wsGastos.Range("G2:G18").Borders.LineStyle = 1
$wsGastos.Columns.Item(7).ColumnWidth = 14.28515625

$wsGastos.Range("A1:G18").Select() | Out-Null

# ---------------------------------------------------------------
# 2) Salas: new sheet right after "Gastos" holding the new table.
# ---------------------------------------------------------------
$wsSalas = $wb.Worksheets.Add($null, $wsGastos)
$wsSalas.Name = "Salas"

$wsSalas.Range("A1").Value = "ID_Sala"
$wsSalas.Range("B1").Value = "NomeSala"
$wsSalas.Range("C1").Value = "ValorTotalGastos"
$wsSalas.Range("A1:C1").Font.Bold = $true
$wsSalas.Range("A1:C1").Font.Color = 0xFFFFFF
$wsSalas.Range("A1:C1").Interior.Color = 0x000000
$wsSalas.Range("A1:C1").HorizontalAlignment = -4108
$wsSalas.Range("A2:C6").Borders.LineStyle = 1

$wsSalas.Columns.Item(1).ColumnWidth = 12
$wsSalas.Columns.Item(2).ColumnWidth = 14.5703125
$wsSalas.Columns.Item(3).ColumnWidth = 20.85546875

$loSalas = $wsSalas.ListObjects.Add(1, $wsSalas.Range("A1:C6"), $null, 1)
$loSalas.Name = "Salas"
$loSalas.TableStyle = "TableStyleDark1"

$wsSalas.Range("A1:C6").Select() | Out-Null

# ---------------------------------------------------------------
# 3) Adicionar_Gastos: reset the ID_Sala seed rows to 1 / 2, wipe
#    the rest of the sample values.
# ---------------------------------------------------------------
$wsAdic = $wb.Worksheets.Item("Adicionar_Gastos")
$wsAdic.Range("B2:E3").ClearContents()
$wsAdic.Range("A2").Value = "1"
$wsAdic.Range("A3").Value = "2"

# ---------------------------------------------------------------
# 4) Home: brand-new, empty sheet at the end; becomes the active tab.
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsHome = $wb.Worksheets.Add($null, $lastSheet)
$wsHome.Name = "Home"
$wsHome.Activate()
